$d = $word.ActiveDocument

# Locate the last paragraph in the document body (the one ending with
# "...în cazul nostru flannel") and append a brand-new paragraph right
# after it, mirroring the same paragraph/run formatting (pStyle Normal,
# spacing before=0/after=160, lang ro-RO).
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastRange = $lastPara.Range

# Insert a new paragraph mark after the last paragraph, then set the
# text of the freshly created paragraph.
$lastRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newRange = $newPara.Range
$newRange.Text = "gata avem script automat de creeare a unui cluster yupy, trb puțin optimizat ca mănâncă multe resurse acum"
